$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 125001176
$ws.Range("J18").Value = 1000
$ws.Range("L18").Value = 1000
$ws.Range("N18").Value = -1568
$ws.Range("H28").Value = 1039.85
$ws.Range("I28").Value = 502.05884
$ws.Range("J28").Value = 4087.3333
$ws.Range("K28").Value = 502.05884
$ws.Range("L28").Value = 4087.3333
$ws.Range("M28").Value = -17.05883999999998
$ws.Range("N28").Value = -5057.3333
$ws.Range("H33").Value = 1632.7142
$ws.Range("J33").Value = 851.3333
$ws.Range("L33").Value = 851.3333
$ws.Range("N33").Value = -1309.3333
$ws.Range("H76").Value = 27783532
$ws.Range("I76").Value = 5099.8335
$ws.Range("K76").Value = 5099.8335
$ws.Range("M76").Value = -4784.8335
$ws.Range("H79").Value = 27783532
$ws.Range("I79").Value = 5099.8335
$ws.Range("K79").Value = 5099.8335
$ws.Range("M79").Value = -4007.8335
$ws.Range("H116").Value = 12509583
$ws.Range("I116").Value = 25006868
$ws.Range("K116").Value = 25006868
$ws.Range("M116").Value = -25003426
$ws.Range("H125").Value = 100000850
$ws.Range("I125").Value = 125000470
$ws.Range("K125").Value = 1125004230
$ws.Range("M125").Value = -1125001770
$ws.Range("H129").Value = 992.3077
$ws.Range("I129").Value = 614.8946999999999
$ws.Range("K129").Value = 1844.6841
$ws.Range("M129").Value = 3155.3159
$ws.Range("H132").Value = 1193.5652
$ws.Range("I132").Value = 1042.75
$ws.Range("K132").Value = 3128.25
$ws.Range("M132").Value = -598.25
$ws.Range("H137").Value = 3637.3125
$ws.Range("I137").Value = 3156.7144
$ws.Range("K137").Value = 9470.143199999999
$ws.Range("M137").Value = -6920.143199999999

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 43479868
$ws.Range("I2").Value = 1324.4
$ws.Range("K2").Value = 1324.4
$ws.Range("M2").Value = -1211.4
$ws.Range("H116").Value = 43479868
$ws.Range("I116").Value = 1324.4
$ws.Range("K116").Value = 1324.4
$ws.Range("M116").Value = 969.5999999999999
$ws.Range("H132").Value = 3783.9253
$ws.Range("I132").Value = 2746.5557
$ws.Range("J132").Value = 8093
$ws.Range("K132").Value = 8239.667099999999
$ws.Range("L132").Value = 24279
$ws.Range("M132").Value = -5709.667099999999
$ws.Range("N132").Value = -29339

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 43479868
$ws.Range("I3").Value = 1324.4
$ws.Range("K3").Value = 1324.4
$ws.Range("M3").Value = -1210.4
$ws.Range("H94").Value = 1481.8959
$ws.Range("I94").Value = 652.5
$ws.Range("K94").Value = 652.5
$ws.Range("M94").Value = -201.5
$ws.Range("H105").Value = 3507.6191
$ws.Range("I105").Value = 2768.0908
$ws.Range("K105").Value = 2768.0908
$ws.Range("M105").Value = -1021.0908

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3488.325
$ws.Range("I16").Value = 1343.5
$ws.Range("K16").Value = 1343.5
$ws.Range("M16").Value = -1056.5
$ws.Range("H94").Value = 1330.9445
$ws.Range("J94").Value = 1309.7
$ws.Range("L94").Value = 1309.7
$ws.Range("N94").Value = -2211.7
$ws.Range("H105").Value = 23811520
$ws.Range("I105").Value = 35714784
$ws.Range("K105").Value = 35714784
$ws.Range("M105").Value = -35713037
$ws.Range("H113").Value = 3488.325
$ws.Range("I113").Value = 1343.5
$ws.Range("K113").Value = 1343.5
$ws.Range("M113").Value = 826.5
$ws.Range("H122").Value = 4168.982
$ws.Range("I122").Value = 3551.45
$ws.Range("J122").Value = 5712.8125
$ws.Range("K122").Value = 10654.35
$ws.Range("L122").Value = 17138.4375
$ws.Range("M122").Value = -8204.349999999999
$ws.Range("N122").Value = -22038.4375
$ws.Range("H132").Value = 3115.1914
$ws.Range("I132").Value = 2263.5588
$ws.Range("K132").Value = 6790.676399999999
$ws.Range("M132").Value = -4260.676399999999
$ws.Range("H134").Value = 3659.9714
$ws.Range("I134").Value = 2119.3845
$ws.Range("K134").Value = 6358.1535
$ws.Range("M134").Value = -3823.1535
$ws.Range("H141").Value = 727842.75
$ws.Range("J141").Value = 727842.75
$ws.Range("L141").Value = 727842.75
$ws.Range("N141").Value = -738202.75

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2552.7144
$ws.Range("I11").Value = 1373.8
$ws.Range("K11").Value = 4121.4
$ws.Range("M11").Value = -3981.4
$ws.Range("H115").Value = 1445.7142
$ws.Range("I115").Value = 1445.7142
$ws.Range("K115").Value = 4337.142599999999
$ws.Range("M115").Value = -3162.142599999999
$ws.Range("H134").Value = 39507.07
$ws.Range("I134").Value = 51032
$ws.Range("K134").Value = 153096
$ws.Range("M134").Value = -148026
$ws.Range("H138").Value = 67129.875
$ws.Range("I138").Value = 67129.875
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 201389.625
$ws.Range("L138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -196249.625

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 20000
$ws.Range("J18").Value = 20000
$ws.Range("L18").Value = 20000
$ws.Range("N18").Value = -20586
$ws.Range("H70").Value = 7370.7617
$ws.Range("I70").Value = 4837.3335
$ws.Range("K70").Value = 4837.3335
$ws.Range("M70").Value = -4567.3335
$ws.Range("H73").Value = 7370.7617
$ws.Range("I73").Value = 4837.3335
$ws.Range("K73").Value = 4837.3335
$ws.Range("M73").Value = -3901.3335
$ws.Range("H113").Value = 6980.722
$ws.Range("I113").Value = 2797.0715
$ws.Range("J113").Value = 9643.046
$ws.Range("K113").Value = 2797.0715
$ws.Range("L113").Value = 9643.046
$ws.Range("M113").Value = -627.0715
$ws.Range("N113").Value = -13983.046
$ws.Range("H126").Value = 20836932
$ws.Range("I126").Value = 50002056
$ws.Range("J126").Value = 4702.643
$ws.Range("K126").Value = 150006168
$ws.Range("L126").Value = 14107.929
$ws.Range("M126").Value = -150003698
$ws.Range("N126").Value = -19047.929
$ws.Range("H132").Value = 2649.8333
$ws.Range("I132").Value = 1529.0834
$ws.Range("K132").Value = 4587.2502
$ws.Range("M132").Value = -2057.2502

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 359.8889
$ws.Range("I16").Value = 359.8889
$ws.Range("K16").Value = 359.8889
$ws.Range("M16").Value = -189.8889
$ws.Range("H22").Value = 2326.625
$ws.Range("J22").Value = 2935.5
$ws.Range("L22").Value = 2935.5
$ws.Range("N22").Value = -3525.5
$ws.Range("H27").Value = 2326.625
$ws.Range("J27").Value = 2935.5
$ws.Range("L27").Value = 2935.5
$ws.Range("N27").Value = -3149.5
$ws.Range("H40").Value = 29415692
$ws.Range("I40").Value = 41668896
$ws.Range("K40").Value = 41668896
$ws.Range("M40").Value = -41668760
$ws.Range("H55").Value = 321.4
$ws.Range("I55").Value = 56
$ws.Range("K55").Value = 56
$ws.Range("M55").Value = 117
$ws.Range("H132").Value = 10212400
$ws.Range("I132").Value = 23812688
$ws.Range("J132").Value = 12183.393
$ws.Range("K132").Value = 71438064
$ws.Range("L132").Value = 36550.179
$ws.Range("M132").Value = -71435534
$ws.Range("N132").Value = -41610.179

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 76927380
$ws.Range("I126").Value = 90911820
$ws.Range("K126").Value = 272735460
$ws.Range("M126").Value = -272732990
$ws.Range("H136").Value = 20615240
$ws.Range("I136").Value = 41667348
$ws.Range("J136").Value = 405214.28
$ws.Range("K136").Value = 125002044
$ws.Range("L136").Value = 1215642.84
$ws.Range("M136").Value = -124999494
$ws.Range("N136").Value = -1220742.84
